$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the full 2x2 cross-join of Sending/Target clusters (ECs, MuSCs) for the Gal-Galr1
# ligand-receptor pair. Rows 2-3 already existed (now updated with TPM-based values);
# rows 4-5 are newly added.
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gal"
$ws.Range("C2").Value = "Galr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03644533333333334
$ws.Range("H2").Value = 0.109336
$ws.Range("I2").Value = 0.005561955322140003
$ws.Range("J2").Value = 0.005561955322140003
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05571033333333333
$ws.Range("N2").Value = 0.167131
$ws.Range("O2").Value = 0.09700855792269546
$ws.Range("P2").Value = 0.09700855792269544
$ws.Range("Q2").Value = 0.002030381668444445
$ws.Range("R2").Value = 0.018273435016
$ws.Range("S2").Value = 0.0005395572650312627
$ws.Range("T2").Value = 0.0005395572650312626
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gal"
$ws.Range("C3").Value = "Galr1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03644533333333334
$ws.Range("H3").Value = 0.109336
$ws.Range("I3").Value = 0.005561955322140003
$ws.Range("J3").Value = 0.005561955322140003
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5185723333333333
$ws.Range("N3").Value = 1.555717
$ws.Range("O3").Value = 0.9029914420773046
$ws.Range("P3").Value = 0.9029914420773046
$ws.Range("Q3").Value = 0.01889954154577778
$ws.Range("R3").Value = 0.170095873912
$ws.Range("S3").Value = 0.00502239805710874
$ws.Range("T3").Value = 0.00502239805710874
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Gal"
$ws.Range("C4").Value = "Galr1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.516166333333334
$ws.Range("H4").Value = 19.548499
$ws.Range("I4").Value = 0.99443804467786
$ws.Range("J4").Value = 0.99443804467786
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05571033333333333
$ws.Range("N4").Value = 0.167131
$ws.Range("O4").Value = 0.09700855792269546
$ws.Range("P4").Value = 0.09700855792269544
$ws.Range("Q4").Value = 0.3630177984854445
$ws.Range("R4").Value = 3.267160186369
$ws.Range("S4").Value = 0.0964690006576642
$ws.Range("T4").Value = 0.09646900065766419
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gal"
$ws.Range("C5").Value = "Galr1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.516166333333334
$ws.Range("H5").Value = 19.548499
$ws.Range("I5").Value = 0.99443804467786
$ws.Range("J5").Value = 0.99443804467786
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5185723333333333
$ws.Range("N5").Value = 1.555717
$ws.Range("O5").Value = 0.9029914420773046
$ws.Range("P5").Value = 0.9029914420773046
$ws.Range("Q5").Value = 3.379103579864778
$ws.Range("R5").Value = 30.411932218783
$ws.Range("S5").Value = 0.897969044020196
$ws.Range("T5").Value = 0.897969044020196
